$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for the Effort-Forecast columns (I/J for CapacityForecast->EffortForecast
# "Forecast" ratio+extrapolation, K/L for CapacityDone->EffortDone "Done" ratio+extrapolation).
$ws.Range("I1").Value = "Forecast"
$ws.Range("K1").Value = "Done"

# Row 2 new formulas.
$ws.Range("I2").Formula = "=IF(D2<>`"`",E2/D2,`"`")"
$ws.Range("J2").Formula = "=I2*D2"
$ws.Range("K2").Formula = "=IF(F2<>`"`",G2/F2,`"`")"
$ws.Range("L2").Formula = "=K2*F2"

# Row 3 gets new CapacityDone / EffortDone values plus the same forecast formulas.
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 7.5

$ws.Range("I3").Formula = "=IF(D3<>`"`",E3/D3,`"`")"
$ws.Range("J3").Formula = "=I3*D3"
$ws.Range("K3").Formula = "=IF(F3<>`"`",G3/F3,`"`")"
$ws.Range("L3").Formula = "=K3*F3"

# Selection as captured in the diff.
$ws.Range("I4:L7").Select()

$wb.Save()
